$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tab-separated: Row, Ticker, Name, Price, MarketCap, Volume, Change(24h)
$data = @"
2	BTC	Bitcoin	36463	713278345221	10638135767	-0.1178
3	ETH	Ethereum	1948.62	234674236247	10594386967	0.28584
4	USDT	Tether	1	87727232727	25595322252	-0.04898
5	BNB	BNB	243.13	37441079374	433375700	-0.19928
6	XRP	XRP	0.611469	32841490138	643918749	1.20883
7	SOL	Solana	59.72	25314152605	2163595384	4.22775
8	USDC	USDC	0.999458	24425990390	5674596136	-0.20277
9	STETH	Lido Staked Ether	1951.33	17527907826	4600282	0.76795
10	ADA	Cardano	0.375245	13141902803	326601194	4.20721
11	DOGE	Dogecoin	0.07828	11137299245	1338360892	-7.97257
12	TRX	TRON	0.102614	9096532280	256923700	0.23388
13	LINK	Chainlink	14.13	7908120394	624270650	5.6966
14	MATIC	Polygon	0.820777	7630221460	439224307	1.11809
15	AVAX	Avalanche	21.25	7584458524	723623568	-0.25366
16	DOT	Polkadot	5.22	6792552843	177863428	1.23394
17	WBTC	Wrapped Bitcoin	36423	5952233544	63126562	0.2034
18	DAI	Dai	1	5290743624	292545455	0.0931
19	LTC	Litecoin	69.04000000000001	5115863614	262479821	-0.31536
20	TON	Toncoin	2.34	5034445988	35889737	1.57073
21	SHIB	Shiba Inu	8.47e-06	5002062801	130075973	-1.59884
22	BCH	Bitcoin Cash	227.02	4445785916	77077961	0.43486
23	UNI	Uniswap	5.06	3827734223	113650291	1.65777
24	LEO	LEO Token	4.08	3786576519	1208205	0.14174
25	OKB	OKB	57.12	3423812857	14897677	1.72499
26	XLM	Stellar	0.118755	3325107382	49276479	0.92264
27	TUSD	TrueUSD	0.998617	3304441994	130823149	-0.12079
28	KAS	Kaspa	0.143476	3113433913	127212618	9.01388
29	XMR	Monero	159.89	2906522015	64871691	-0.47624
30	ETC	Ethereum Classic	19.03	2732624860	97990304	0.57941
31	ATOM	Cosmos Hub	9.109999999999999	2668287560	184072785	-0.55987
32	CRO	Cronos	0.096206	2538933153	17683434	-1.19075
33	FIL	Filecoin	4.73	2232087839	154823470	3.01535
34	HBAR	Hedera	0.060779	2045255233	44063818	-0.46653
35	LDO	Lido DAO	2.27	2025086884	94278091	4.31648
36	ICP	Internet Computer	4.43	1999639978	57417254	6.10653
37	APT	Aptos	6.97	1927349159	90417905	0.76695
38	NEAR	NEAR Protocol	1.81	1819106245	220442547	4.21148
39	BUSD	BUSD	0.999793	1803129302	2178141636	-0.09601999999999999
40	RUNE	THORChain	5.37	1626569223	606647457	-13.20237
41	IMX	Immutable	1.29	1610125084	746641476	19.17595
42	MNT	Mantle	0.496844	1546042804	16626316	-1.81646
43	OP	Optimism	1.72	1512620694	110625115	2.44179
44	VET	VeChain	0.02075484	1511557383	29933293	0.17407
45	TAO	Bittensor	256.86	1460659872	20209119	19.43234
46	QNT	Quant	98.78	1437243683	19494840	-0.11649
47	INJ	Injective	15.68	1317429856	68947999	-0.54709
48	ARB	Arbitrum	1.02	1303115892	227473803	-0.44736
49	AAVE	Aave	88.33	1296056485	130091736	2.99954
50	RNDR	Render	3.42	1285194729	247908508	9.540710000000001
51	GRT	The Graph	0.131363	1226772951	53392678	5.22045
"@

$lines = $data -split "`r?`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split "`t"
    $r = [int]$parts[0]
    $ws.Cells.Item($r, 2).Value = $parts[1]
    $ws.Cells.Item($r, 3).Value = $parts[2]
    $ws.Cells.Item($r, 4).Value = [double]$parts[3]
    $ws.Cells.Item($r, 5).Value = [double]$parts[4]
    $ws.Cells.Item($r, 6).Value = [double]$parts[5]
    $ws.Cells.Item($r, 7).Value = [double]$parts[6]
}
